$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data base")

# --- Resource (column C): citrus host plants for the four Japan rows ---
$ws.Range("C18").Value = "Citrus unshiu"
$ws.Range("C19").Value = "Citrus aurantium"
$ws.Range("C20").Value = "Citrus unshiu"
$ws.Range("C21").Value = "Citrus unshiu"

# --- Used (column G): updated fitted estimates of r --------------------
$ws.Range("G18").Value = "No"
$ws.Range("G19").Value = "Yes"
$ws.Range("G20").Value = "Yes"
$ws.Range("G21").Value = "No"

# --- Notes (column K) ----------------------------------------------------
$ws.Range("K18").Value = "Not used because only 4 temperatures for most traits"
$ws.Range("K18").Font.Size = 11

$ws.Range("K21").Value = "Not used because development is not unimodal"
$ws.Range("K21").Font.Size = 11

# --- Restore selection / scroll position --------------------------------
$ws.Activate() | Out-Null
$ws.Range("B28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
